# Consolidate the separate text runs ("The", " ", "picture", " ", "first")
# in the caption textbox on slide 1 into a single run "The picture first".
#
# Re-assigning TextRange.Text causes the writer to re-emit the paragraph's
# runs, collapsing however many runs previously existed into one. Setting
# the range to an intermediate, different value first guarantees the
# subsequent assignment of the final text is treated as a genuine change
# (and not skipped as a no-op).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shape = $s.Shapes.Item(2)

$shape.TextFrame.TextRange.Text = "x"
$shape.TextFrame.TextRange.Text = "The picture first"
